$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Cluster Name"

# Data rows: new cluster names (re-sorted alphabetically) and updated active-case counts
$data = @(
    @{Row=2; Name='3035 Campbell Place Aged Care Glen Waverley'; Value=10}
    @{Row=3; Name='3622 Olivet Care Aged Care Services Ringwood'; Value=12}
    @{Row=4; Name='3961 Heritage Care Water Gardens Aged CareFacility Sydenham'; Value=26}
    @{Row=5; Name='4167 Royal Freemasons Centennial LodgeWantirna South'; Value=20}
    @{Row=6; Name='AG Industries Pty Ltd Factory Thomastown'; Value=17}
    @{Row=7; Name='Aintree Primary School Aintree'; Value=17}
    @{Row=8; Name='Assisi Centre Aged Care Rosanna'; Value=19}
    @{Row=9; Name='Australian Meat Group Abattoir DandenongSouth'; Value=13}
    @{Row=10; Name='Bacchus Marsh Childcare and KindergartenCentre Bacchus Marsh'; Value=34}
    @{Row=11; Name='Baden Powell College Tarneit'; Value=15}
    @{Row=12; Name='Bandiana Primary School Bandiana'; Value=10}
    @{Row=13; Name='Covenant College Bell Post Hill'; Value=17}
    @{Row=14; Name='Domestic Freight Vessel Siem AquamarineGeelong'; Value=10}
    @{Row=15; Name='Gladstone Park Secondary College'; Value=14}
    @{Row=16; Name='Guardian Childcare & Education Moorabbin'; Value=11}
    @{Row=17; Name='Hamlyn Banks Primary School Hamlyn Heights'; Value=11}
    @{Row=18; Name='Hazelwood North Primary School HazelwoodNorth'; Value=27}
    @{Row=19; Name='Islamic College of Melbourne Tarneit'; Value=23}
    @{Row=20; Name='Master Poultry Group West Footscray'; Value=13}
    @{Row=21; Name='Morwell Park Primary School Morwell'; Value=58}
    @{Row=22; Name='Nido Early School Woodend'; Value=11}
    @{Row=23; Name='Northern Bay College Goldsworthy 9-12 CampusCorio'; Value=18}
    @{Row=24; Name='Northern Bay College Wexford Campus Corio'; Value=40}
    @{Row=25; Name='Northern Health Northern Hospital EppingEmergency Department'; Value=21}
    @{Row=26; Name='Northern Health The Northern Hospital Epping'; Value=12}
    @{Row=27; Name='Oakleigh South Primary School Oakleigh South'; Value=16}
    @{Row=28; Name='Our Lady''s Catholic Primary School Wangaratta'; Value=12}
    @{Row=29; Name='Rutherglen Motor Inn and Walkabout MotelRutherglen'; Value=18}
    @{Row=30; Name='Sirius College Ibrahim Dellal Campus Sunshine'; Value=11}
    @{Row=31; Name='Sirius College Shepparton Campus Shepparton'; Value=11}
    @{Row=32; Name='Smartie Pants Early Learning and DevelopmentDiamond Creek'; Value=19}
    @{Row=33; Name='St Brendans Primary School Shepparton'; Value=10}
    @{Row=34; Name='St Clare''s Primary School Officer'; Value=10}
    @{Row=35; Name='St Georges Road Primary School Shepparton'; Value=13}
    @{Row=36; Name='St Joseph''s School Quarry Hill'; Value=31}
    @{Row=37; Name='St Josephs Catholic Primary School Warragul'; Value=12}
    @{Row=38; Name='St Louis de Montfort''s School Aspendale'; Value=13}
    @{Row=39; Name='St Paul''s Primary School Sunshine West'; Value=12}
    @{Row=40; Name='St Vincents Hospital Emergency DepartmentMelbourne'; Value=12}
    @{Row=41; Name='Stockdale Road Primary School Traralgon'; Value=33}
    @{Row=42; Name='Story House Early Learning Epping'; Value=12}
    @{Row=43; Name='Sunbury Primary School Sunbury'; Value=11}
    @{Row=44; Name='TUROSI PTY LTD Thomastown'; Value=11}
    @{Row=45; Name='Templestowe Park Primary School Templestowe'; Value=29}
    @{Row=46; Name='The Lake Primary School Cabarita'; Value=19}
    @{Row=47; Name='Warragul Regional College Warragul'; Value=11}
    @{Row=48; Name='Werribee Mercy Hospital Emergency Department'; Value=21}
    @{Row=49; Name='Western Health Sunshine Hospital Emergency Department St Albans'; Value=14}
    @{Row=50; Name='Wodonga Primary School Wodonga'; Value=15}
    @{Row=51; Name='Wodonga Senior Secondary College Wodonga'; Value=13}
    @{Row=52; Name='Wodonga South Primary School Wodonga'; Value=15}
    @{Row=53; Name='Wyndham Christian College Wyndham Vale'; Value=14}
    @{Row=54; Name='Yeshivah College St Kilda East'; Value=10}
)

foreach ($row in $data) {
    $ws.Cells.Item($row.Row, 1).Value = $row.Name
    $ws.Cells.Item($row.Row, 2).Value = $row.Value
}
